$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: update Price (D) and Volume(1h) (E) columns.
# A handful of new Price strings (e.g. "1.00", "2.50") are valid numeric
# literals; format those cells as Text first so Excel keeps the exact
# scraped string instead of collapsing it to a float (1 / 2.5) or losing
# the trailing zero.

$ws.Range("D2").Value = "69.765.35"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "3.757.20"
$ws.Range("E3").Value = "  +2.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.66"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.29"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").Value = "3.755.31"
$ws.Range("E7").Value = "  +2.49%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  +3.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").Value = "  -4.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("E12").Value = "  -1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.24"
$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").Value = "4.375.79"
$ws.Range("E15").Value = "  +2.43%  "

$ws.Range("D16").Value = "3.764.57"
$ws.Range("E16").Value = "  +2.82%  "

$ws.Range("D17").Value = "69.886.43"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.123"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.74"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "508.04"
$ws.Range("E21").Value = "  -2.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.57"
$ws.Range("E22").Value = "  +3.06%  "

$ws.Range("E23").Value = "  -2.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.27"
$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.16"
$ws.Range("E26").Value = "  -2.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.15"
$ws.Range("E27").Value = "  +1.79%  "

$ws.Range("E28").Value = "  +24.77%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.89"
$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.14"
$ws.Range("E33").Value = "  -1.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  +4.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.22"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("E39").Value = "  +2.56%  "

$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.21"
$ws.Range("E41").Value = "  -2.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.86"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "428.71"
$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.74"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("E45").Value = "  +2.98%  "

$ws.Range("D46").Value = "3.006.54"
$ws.Range("E46").Value = "  -3.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("E47").Value = "  -1.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.43"
$ws.Range("E48").Value = "  -4.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.15"
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.50"
$ws.Range("E51").Value = "  +1.36%  "
